$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","H","J","L","O")

$data = @(
    @(19.93587785526115,9.474960569113012,8.10468717712231,12.69313455404209,35.67475656164122,7.344005520526261,9.931293539833714,11.21381394625839,27.43644700732835),
    @(19.39582765267504,9.18825208512923,8.098123691196937,12.73006694807218,35.81446367997784,7.344005520526261,9.95982585977629,11.1838795834798,27.56367958688883),
    @(19.05820643016829,9.006950416796826,8.094827478781447,12.75427584390294,35.90987989271635,7.344005520526261,9.978246839644276,11.16662915140331,27.6488708392633),
    @(18.91929841816443,8.931830301089015,8.093669891693935,12.76452704375945,35.95117724214588,7.344005520526261,9.985981044287067,11.15988749693787,27.68535949497018),
    @(18.89615866105435,8.919284499285537,8.093488919649712,12.76625257090615,35.95818021146511,7.344005520526261,9.987279064026286,11.15878556121462,27.69152531028196),
    @(19.05633817026228,9.00594221670935,8.094811114026127,12.75441253198799,35.91042707852613,7.344005520526261,9.978350223763064,11.1665370596509,27.64935576815508),
    @(19.75103075546568,9.377244163334673,8.102272479336522,12.70555117263805,35.72092272928676,7.344005520526261,9.940944683314559,11.20326065523002,27.47884628519506),
    @(21.05770306077975,10.06025832788164,8.122675681328774,12.62186580809286,35.42612554349291,7.344005520526261,9.874718186583646,11.28405212049552,27.20080800569252),
    @(21.9742759935203,10.53056907621636,8.141114242718215,12.56774175984621,35.2568551386231,7.344005520526261,9.830361249457708,11.34850287659246,27.03119049047471),
    @(22.38005833456714,10.7369749017233,8.150235116588949,12.54470980996103,35.19022169840525,7.344005520526261,9.81110631782594,11.3788706241971,26.96161960103012),
    @(22.53198171655912,10.81400007320526,8.153792882682161,12.53621619690632,35.16648782045997,7.344005520526261,9.80394703000907,11.39051563926581,26.93637170920299),
    @(22.49934153505563,10.79746262307897,8.153022058547746,12.53803531149341,35.17153255836219,7.344005520526261,9.805483043843909,11.38800129761143,26.94176040920655),
    @(22.39259269388031,10.74333485725722,8.150525743010274,12.54400646672051,35.18823901746553,7.344005520526261,9.810514674262201,11.37982578392728,26.9595204213738),
    @(22.32697597961223,10.71003062635889,8.149010160763741,12.54769366177008,35.19866761290062,7.344005520526261,9.813613882779137,11.37483681901116,26.97054198677957),
    @(21.94752087255055,10.51692366078782,8.140532794442258,12.56927889423806,35.2614191057509,7.344005520526261,9.831638129198089,11.34653891864663,27.03589031951618),
    @(21.71177894044969,10.39648861338658,8.13551889672736,12.58292747640889,35.30257639969693,7.344005520526261,9.842931442350643,11.32944389687432,27.07792727926874),
    @(21.57514195958138,10.32651048857603,8.132704131509596,12.59092740419662,35.32722447922629,7.344005520526261,9.849513996285649,11.3197103423102,27.10281979325784),
    @(21.52870383960415,10.30269733443693,8.131763013234828,12.59366175455112,35.33573719381209,7.344005520526261,9.851757686974913,11.31643190186329,27.11137041704164),
    @(21.73698315834112,10.40938270318708,8.136045494387588,12.58145907933873,35.29809411801483,7.344005520526261,9.84172025658863,11.33125347847752,27.07337843791521),
    @(22.42399555788075,10.75926471585162,8.151256164338303,12.54224640711762,35.18329118952657,7.344005520526261,9.809033181254541,11.38222322800762,26.95427405455437),
    @(22.86281765702444,10.98128993499753,8.161802011738683,12.51794791324244,35.11700008018965,7.344005520526261,9.788440195462902,11.41637983770337,26.88282979775444),
    @(22.62958170550695,10.86341437870565,8.156118675392408,12.5307950000777,35.15157878585806,7.344005520526261,9.799360818851314,11.39807432735963,26.92037368598054),
    @(21.72559176566832,10.40355558122439,8.135807208323806,12.58212246448951,35.30011748613899,7.344005520526261,9.842267553654169,11.33043507156155,27.07543271143609),
    @(20.71120230853495,9.880753759855075,8.116545193993081,12.64320990978359,35.49760049959049,7.344005520526261,9.891875957499948,11.28405212049552,27.26996128914769)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$rowNum").Value = $rowVals[$j]
    }
}
